# Hyperwall_stats.xlsx update - "new version of coverage figures"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet2: header row additions
# ---------------------------------------------------------------------------
$ws2.Range("C1").Value = "num_files"
$ws2.Range("E1").Value = "in GB"
$ws2.Range("H1").Value = "from SR15 (GB)"

# ---------------------------------------------------------------------------
# Sheet2: new "num_files" (C) column, rows 2-21
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value = 2725
$ws2.Range("C3").Value = 9842
$ws2.Range("C4").Value = 1532075
$ws2.Range("C5").Value = 18840050
$ws2.Range("C6").Value = 2248346
$ws2.Range("C7").Value = 9406271
$ws2.Range("C8").Value = 17422033
$ws2.Range("C9").Value = 8726
$ws2.Range("C10").Value = 955434
$ws2.Range("C12").Value = 3342905
$ws2.Range("C13").Value = 1596
$ws2.Range("C14").Value = 6728700
$ws2.Range("C15").Value = 12408
$ws2.Range("C16").Value = 104280276
$ws2.Range("C18").Value = 1407298
$ws2.Range("C19").Value = 1353640
$ws2.Range("C20").Value = 4678
$ws2.Range("C21").Value = 1429

# ---------------------------------------------------------------------------
# Sheet2: "sum_contentLength" (D) column - new values / updated values (bytes)
# ---------------------------------------------------------------------------
$ws2.Range("D2").Value = 2451198076
$ws2.Range("D3").Value = 81583137856
$ws2.Range("D4").Value = 587677729967
$ws2.Range("D5").Value = 17990760357811
$ws2.Range("D6").Value = 57908050171200
$ws2.Range("D7").Value = 9379144910661
$ws2.Range("D8").Value = 146618669070408
$ws2.Range("D9").Value = 3143178514
$ws2.Range("D10").Value = 191035141255
$ws2.Range("D12").Value = 4403495417547
$ws2.Range("D13").Value = 216871527936
$ws2.Range("D14").Value = 17350990159044
$ws2.Range("D15").Value = 1687020315648
$ws2.Range("D16").Value = 1177090480293440
$ws2.Range("D18").Value = 7621203911483
$ws2.Range("D19").Value = 34821786831339
$ws2.Range("D20").Value = 218829672
$ws2.Range("D21").Value = 15079956

# rows with no byte total yet (JWST, SPITZER_SHA) - leave the cell present but empty
$ws2.Range("D11").NumberFormat = "General"
$ws2.Range("D17").NumberFormat = "General"

# ---------------------------------------------------------------------------
# Sheet2: helper constants (1024^3 and 1 billion) used by the "in GB" column
# ---------------------------------------------------------------------------
$ws2.Range("D24").Value = "1024^3"
$ws2.Range("E24").Formula = '=1024*1024*1024'

$ws2.Range("D25").Value = "1 billion"
$ws2.Range("E25").Value = 1000000000
$ws2.Range("D25:E25").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Sheet2: "in GB" (E) column formulas, rows 2-21
# ---------------------------------------------------------------------------
$ws2.Range("E2").Formula = '=D2/$E$24'
$ws2.Range("E3:E21").Formula = '=D3/$E$24'

$ws2.Range("E2:E19").NumberFormat = "_(* #,##0.0_);_(* \(#,##0.0\);_(* ""-""??_);_(@_)"
$ws2.Range("E20:E21").NumberFormat = "_(* #,##0.000_);_(* \(#,##0.000\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# Sheet2: "from SR15 (GB)" (H) column, a handful of rows
# ---------------------------------------------------------------------------
$ws2.Range("H3").Value = 96
$ws2.Range("H4").Value = 1200
$ws2.Range("H5").Value = 28590
$ws2.Range("H6").Value = 90000
$ws2.Range("H7").Value = 7505
$ws2.Range("H8").Value = 107230
$ws2.Range("H10").Value = 600
$ws2.Range("H18").Value = 6641
$ws2.Range("H4:H8").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Sheet2: move the two small helper sums down below the new rows
# ---------------------------------------------------------------------------
$ws2.Range("G3").ClearContents()
$ws2.Range("G26").Value = "shuttle"
$ws2.Range("G27").Formula = '=SUM(B2,B9,B20:B21)'
$ws2.Range("H27").Formula = '=E2+E9+E20+E21'
$ws2.Range("H27").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# Sheet2: totals row (replaces the old lone SUM(B2:B21) in B23)
# ---------------------------------------------------------------------------
$ws2.Range("B23").ClearContents()

$ws2.Range("A28").Value = "totals"
$ws2.Range("B28").Formula = '=SUM(B2:B21)'
$ws2.Range("C28").Formula = '=SUM(C2:C21)'
$ws2.Range("E28").Formula = '=SUM(E2:E21)'
$ws2.Range("B28:C28").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws2.Range("E28").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$ws2.Range("D28").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# Sheet2: column widths for the newly-populated columns
# ---------------------------------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 14.1666666666667
$ws2.Columns.Item(3).ColumnWidth = 13.1666666666667
$ws2.Columns.Item(5).ColumnWidth = 12.6666666666667

# ---------------------------------------------------------------------------
# Selections: set Sheet1's selection first, then Sheet2's last so Sheet2
# remains the active / displayed tab (matching tabSelected + activeTab).
# ---------------------------------------------------------------------------
$ws1.Range("F25").Select()
$ws2.Range("I28").Select()
